$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Tabelle1: new "Mage Guild" ability block (Scalebreaker skill detection) ---
$ws1.Range("A134").Value = "Mage Guild"

$ws1.Range("A136").Value = "Entropy"
$ws1.Range("B136").Value = 28567
$ws1.Range("C136").Value = 126370
$ws1.Range("D136").Value = "Entropy"
$ws1.Range("E136").Value = 2240

$ws1.Range("A137").Value = "Degeneration"
$ws1.Range("B137").Value = 40457
$ws1.Range("C137").Value = 126374
$ws1.Range("D137").Value = "Degeneration"
$ws1.Range("E137").Value = 2240

$ws1.Range("A138").Value = "Structured Entropy"
$ws1.Range("B138").Value = 40452
$ws1.Range("C138").Value = 126371
$ws1.Range("D138").Value = "Structured Entropy"
$ws1.Range("E138").Value = 2240

$ws1.Range("A140").Value = "Soul Magic"

$ws1.Range("A142").Value = "Soul Trap"
$ws1.Range("B142").Value = 26768
$ws1.Range("C142").Value = 126891
$ws1.Range("D142").Value = "Soul Trap"
$ws1.Range("E142").Value = 2240

$ws1.Range("A143").Value = "Soul Splitting Trap"
$ws1.Range("B143").Value = 40328
$ws1.Range("C143").Value = 126894
$ws1.Range("D143").Value = "Soul Splitting Trap"
$ws1.Range("E143").Value = 2240

$ws1.Range("A144").Value = "Consuming Trap"
$ws1.Range("B144").Value = 40317
$ws1.Range("C144").Value = 126896
$ws1.Range("D144").Value = "Consuming Trap"
$ws1.Range("E144").Value = 2240

# --- Tabelle2: extend the lookup/report formulas down to match the new Tabelle1 rows ---
for ($n = 135; $n -le 147; $n++) {
    $m = $n + 3
    $ws2.Range("A$n").Formula = "=IF(ISBLANK(Tabelle1!A$m),`"`",Tabelle1!A$m)"
    $ws2.Range("B$n").Formula = "=IF(ISBLANK(Tabelle1!B$m),`"nil`",Tabelle1!B$m)"
    $ws2.Range("C$n").Formula = "=IF(ISBLANK(Tabelle1!C$m),`"nil`",Tabelle1!C$m)"
    $ws2.Range("D$n").Formula = "=IF(ISBLANK(Tabelle1!D$m),`"`",Tabelle1!D$m)"
    $ws2.Range("E$n").Formula = "=IF(ISBLANK(Tabelle1!E$m),`"nil`",Tabelle1!E$m)"
    $ws2.Range("F$n").Formula = "=IF(ISBLANK(Tabelle1!F$m),`"nil`",Tabelle1!F$m)"
    $ws2.Range("G$n").Formula = "=IF(ISBLANK(Tabelle1!G$m),`"nil`",Tabelle1!G$m)"
    $ws2.Range("I$n").Formula = "=IF(ISBLANK(Tabelle1!B$m),`"`",CONCATENATE(`"[`",B$n,`"] = {`",_xlfn.TEXTJOIN(`", `", FALSE, Tabelle2!C$n,Tabelle2!E${n}:G$n),`"}, --`", A$n, `" --> `",D$n))"
}

# The existing Tabelle2 rows 131-134 (which already referenced Tabelle1 rows
# 134-137 before this edit) have their "I" column formulas cross-referencing
# other cells on Tabelle2 itself; re-apply those formulas so their cached
# values pick up the newly entered Tabelle1 data.
for ($n = 131; $n -le 134; $n++) {
    $ws2.Range("I$n").Formula = $ws2.Range("I$n").Formula
}

# --- View / window state ---
# Tabelle1 becomes the active sheet, with the header row frozen and the
# viewport scrolled down near the newly added rows.
$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 113
$ws1.Range("D138").Select()

# Tabelle2 keeps its own scroll/selection state, but is no longer the active tab.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 116
$ws2.Range("I140").Select()

$ws1.Activate()

# --- Page setup for Tabelle1 (paper size / orientation) ---
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

$wb.Application.CalculateFull()
